$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(234, 44308, 1, 8, 151.5151515151515),
    @(235, 44309, 1, 7, 132.5757575757576),
    @(236, 44310, 1, 8, 151.5151515151515),
    @(237, 44311, 1, 6, 113.6363636363636),
    @(238, 44312, 0, 4, 75.75757575757575)
)

foreach ($vals in $data) {
    $row = $vals[0]

    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $vals[1]
    $ws.Cells.Item($row, 2).Value = $vals[2]
    $ws.Cells.Item($row, 3).Value = $vals[3]
    $ws.Cells.Item($row, 4).Value = $vals[4]
}

$excel.CutCopyMode = 0
